$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing shared-string cells (rows 3 & 4) with appended test cases ---
$ws.Range("D3").Value = "1. CIR-W0006`n2. CIR-W0007`n3. CIR-W0008`n4. CIR-W0009`n5. CIR-W0012`n6. CIR-W0016`n7. CIR-W0017`n8. CIR-W0018"
$ws.Range("E3").Value = "1. Работа «Зритель-контраст» в «АРМ ОПК» `n2. Захват кадра в главной форме «АРМ ОПК»`n3. Захват кадра на границе окна «Видеокамера»`n4. Отработка таймаута в АРМ ОПК`n5. Вызов формы захвата во время захвата в АРМ ОПК`n6. Захват с одновременным сканированием документа с ИС`n7. Захват с одновременным сканированием документа без ИС`n8. Отсутствие функции удаления фотографии при нажатии по кнопке «Изменить» "
$ws.Range("D4").Value = "1. CIR-W0010`n2. CIR-W0011`n3. CIR-W0013`n4. CIR-W0014`n5. CIR-W0015"
$ws.Range("E4").Value = "1. Вызов формы захвата с наличием захваченного кадра`n2. Вызов формы захвата без захваченного кадра`n3. Нажатие по кнопке «Старт F5» в форме захвата`n4. Нажатие по кнопке «Стоп F6» в форме захвата`n5. Нажатие по кнопке «Закрыть Esc» в форме захвата"

# --- Row 5: new section CIR-W S2.3 ---
$ws.Range("C5").Value = "Проверить отсутствие графических элементов"
$ws.Range("D5").Value = "1. CIR-W0019`n2. CIR-W0020"
$ws.Range("E5").Value = "1. Отсутствие в списке «Настройка» пункта «Зритель-Каскад»`n2. Отсутствие демо панели в главной форме "

# --- Row 6: new section CIR-W S2.4 ---
$ws.Range("C6").Value = "Проверить работу горячих клавиш в АРМ ОПК"
$ws.Range("D6").Value = "1. CIR-W0021`n2. CIR-W0022`n3. CIR-W0023"
$ws.Range("E6").Value = "1. Отмена выбора кадра сочетанием клавиш Ctrl + F7`n2. Открытие формы захвата сочетанием клавиш Ctrl + F8`n3. Запуск захвата сочетанием клавиш Ctrl + F9"

# --- Row 7: new section CIR-W S2.5 ---
$ws.Range("B7").Value = "CIR-W S2.5 "
$ws.Range("C7").Value = "Отключение видоекамеры"
$ws.Range("D7").Value = "1. CIR-W0024"
$ws.Range("E7").Value = "1. Отключение видеокамеры во время захвата "

# --- B5 / B6: set value AND change style from 3 -> 1 (no border, matches header-id style) ---
$ws.Range("B5").Value = "CIR-W S2.3 "
$ws.Range("B5").Font.Name = "Times New Roman"
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("B5").Borders.LineStyle = 0
$ws.Range("B6").Value = "CIR-W S2.4 "
$ws.Range("B6").Font.Name = "Times New Roman"
$ws.Range("B6").Font.Size = 12
$ws.Range("B6").WrapText = $true
$ws.Range("B6").HorizontalAlignment = -4131
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6").Borders.LineStyle = 0

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 23.325
$ws.Columns.Item(5).ColumnWidth = 45.68

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 114.75
$ws.Rows.Item(3).RowHeight = 207.75
$ws.Rows.Item(4).RowHeight = 157.5
$ws.Rows.Item(5).RowHeight = 63
$ws.Rows.Item(6).RowHeight = 81
$ws.Rows.Item(7).RowHeight = 31.5

# --- Selection / view ---
$win = $excel.ActiveWindow
try { $win.ScrollRow = 4 } catch {}
try { $win.ScrollColumn = 2 } catch {}
try { $win.TopLeftCell = $ws.Range("B4") } catch {}
$ws.Range("E7").Select()
